$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (2020-05-20) to the "Condicion_Pacientes" table
$tbl = $ws.ListObjects.Item("Condicion_Pacientes")
$newRow = $tbl.ListRows.Add()

# Inherit the number formatting of the previous data row (68)
$ws.Range("A68:F68").Copy()
$ws.Range("A69:F69").PasteSpecial(-4122)

$ws.Range("A69").Value = 43971
$ws.Range("B69").Value = 540
$ws.Range("C69").Value = 145
$ws.Range("D69").Value = 307
$ws.Range("E69").Value = 13
$ws.Range("F69").Value = 16

$ws.Range("F69").Select()
